$d = $word.ActiveDocument

$d.Content.Find.Execute("24+30=54", $true, $false, $false, $false, $false, $true, 1, $false, "13-1=12", 2) | Out-Null
$d.Content.Find.Execute("88-21=67", $true, $false, $false, $false, $false, $true, 1, $false, "68-47=21", 2) | Out-Null
$d.Content.Find.Execute("42-1=41", $true, $false, $false, $false, $false, $true, 1, $false, "48+10=58", 2) | Out-Null
$d.Content.Find.Execute("68-34=34", $true, $false, $false, $false, $false, $true, 1, $false, "90-15=75", 2) | Out-Null
$d.Content.Find.Execute("20-16=4", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=57", 2) | Out-Null
$d.Content.Find.Execute("29-5=24", $true, $false, $false, $false, $false, $true, 1, $false, "75-13=62", 2) | Out-Null
$d.Content.Find.Execute("19+30=49", $true, $false, $false, $false, $false, $true, 1, $false, "14-0=14", 2) | Out-Null
$d.Content.Find.Execute("36-7=29", $true, $false, $false, $false, $false, $true, 1, $false, "71-21=50", 2) | Out-Null
$d.Content.Find.Execute("20+40=60", $true, $false, $false, $false, $false, $true, 1, $false, "70-62=8", 2) | Out-Null
$d.Content.Find.Execute("63-36=27", $true, $false, $false, $false, $false, $true, 1, $false, "18+17=35", 2) | Out-Null
$d.Content.Find.Execute("67+0=67", $true, $false, $false, $false, $false, $true, 1, $false, "29+17=46", 2) | Out-Null
$d.Content.Find.Execute("50+5=55", $true, $false, $false, $false, $false, $true, 1, $false, "50+12=62", 2) | Out-Null
$d.Content.Find.Execute("99-49=50", $true, $false, $false, $false, $false, $true, 1, $false, "41-13=28", 2) | Out-Null
$d.Content.Find.Execute("35+38=73", $true, $false, $false, $false, $false, $true, 1, $false, "21+11=32", 2) | Out-Null
$d.Content.Find.Execute("78-32=46", $true, $false, $false, $false, $false, $true, 1, $false, "94-30=64", 2) | Out-Null
$d.Content.Find.Execute("23+46=69", $true, $false, $false, $false, $false, $true, 1, $false, "51+47=98", 2) | Out-Null
$d.Content.Find.Execute("2+23=25", $true, $false, $false, $false, $false, $true, 1, $false, "56-23=33", 2) | Out-Null
$d.Content.Find.Execute("50-13=37", $true, $false, $false, $false, $false, $true, 1, $false, "65+3=68", 2) | Out-Null
$d.Content.Find.Execute("42-38=4", $true, $false, $false, $false, $false, $true, 1, $false, "46+39=85", 2) | Out-Null
$d.Content.Find.Execute("97-40=57", $true, $false, $false, $false, $false, $true, 1, $false, "11+12=23", 2) | Out-Null
$d.Content.Find.Execute("3+82=85", $true, $false, $false, $false, $false, $true, 1, $false, "36-17=19", 2) | Out-Null
$d.Content.Find.Execute("88-16=72", $true, $false, $false, $false, $false, $true, 1, $false, "19+59=78", 2) | Out-Null
$d.Content.Find.Execute("78-22=56", $true, $false, $false, $false, $false, $true, 1, $false, "29+33=62", 2) | Out-Null
$d.Content.Find.Execute("27+64=91", $true, $false, $false, $false, $false, $true, 1, $false, "97-33=64", 2) | Out-Null
$d.Content.Find.Execute("68+14=82", $true, $false, $false, $false, $false, $true, 1, $false, "87+0=87", 2) | Out-Null
$d.Content.Find.Execute("58-26=32", $true, $false, $false, $false, $false, $true, 1, $false, "29+34=63", 2) | Out-Null
$d.Content.Find.Execute("39+13=52", $true, $false, $false, $false, $false, $true, 1, $false, "36-33=3", 2) | Out-Null
$d.Content.Find.Execute("37-35=2", $true, $false, $false, $false, $false, $true, 1, $false, "24+43=67", 2) | Out-Null
$d.Content.Find.Execute("76-76=0", $true, $false, $false, $false, $false, $true, 1, $false, "21+62=83", 2) | Out-Null
$d.Content.Find.Execute("26+2=28", $true, $false, $false, $false, $false, $true, 1, $false, "69-43=26", 2) | Out-Null
$d.Content.Find.Execute("38+61=99", $true, $false, $false, $false, $false, $true, 1, $false, "67-34=33", 2) | Out-Null
$d.Content.Find.Execute("18+38=56", $true, $false, $false, $false, $false, $true, 1, $false, "24+73=97", 2) | Out-Null
$d.Content.Find.Execute("46-32=14", $true, $false, $false, $false, $false, $true, 1, $false, "72-30=42", 2) | Out-Null
$d.Content.Find.Execute("47-38=9", $true, $false, $false, $false, $false, $true, 1, $false, "36-3=33", 2) | Out-Null
$d.Content.Find.Execute("72-59=13", $true, $false, $false, $false, $false, $true, 1, $false, "33+60=93", 2) | Out-Null
$d.Content.Find.Execute("52+12=64", $true, $false, $false, $false, $false, $true, 1, $false, "75+9=84", 2) | Out-Null
$d.Content.Find.Execute("89-47=42", $true, $false, $false, $false, $false, $true, 1, $false, "44+12=56", 2) | Out-Null
$d.Content.Find.Execute("93-6=87", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=61", 2) | Out-Null
$d.Content.Find.Execute("51+3=54", $true, $false, $false, $false, $false, $true, 1, $false, "62-41=21", 2) | Out-Null
$d.Content.Find.Execute("89-71=18", $true, $false, $false, $false, $false, $true, 1, $false, "43+53=96", 2) | Out-Null
$d.Content.Find.Execute("16+32=48", $true, $false, $false, $false, $false, $true, 1, $false, "77-26=51", 2) | Out-Null
$d.Content.Find.Execute("30+54=84", $true, $false, $false, $false, $false, $true, 1, $false, "97-20=77", 2) | Out-Null
$d.Content.Find.Execute("29+38=67", $true, $false, $false, $false, $false, $true, 1, $false, "17+28=45", 2) | Out-Null
$d.Content.Find.Execute("49-35=14", $true, $false, $false, $false, $false, $true, 1, $false, "93-37=56", 2) | Out-Null
$d.Content.Find.Execute("25+61=86", $true, $false, $false, $false, $false, $true, 1, $false, "37+18=55", 2) | Out-Null
$d.Content.Find.Execute("63-45=18", $true, $false, $false, $false, $false, $true, 1, $false, "35+12=47", 2) | Out-Null
$d.Content.Find.Execute("53+12=65", $true, $false, $false, $false, $false, $true, 1, $false, "17+13=30", 2) | Out-Null
$d.Content.Find.Execute("20+49=69", $true, $false, $false, $false, $false, $true, 1, $false, "55-24=31", 2) | Out-Null
$d.Content.Find.Execute("61-8=53", $true, $false, $false, $false, $false, $true, 1, $false, "6+14=20", 2) | Out-Null
$d.Content.Find.Execute("40-10=30", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=93", 2) | Out-Null
$d.Content.Find.Execute("34-11=23", $true, $false, $false, $false, $false, $true, 1, $false, "64-33=31", 2) | Out-Null
$d.Content.Find.Execute("68-25=43", $true, $false, $false, $false, $false, $true, 1, $false, "67-51=16", 2) | Out-Null
$d.Content.Find.Execute("4+25=29", $true, $false, $false, $false, $false, $true, 1, $false, "49-8=41", 2) | Out-Null
$d.Content.Find.Execute("26+49=75", $true, $false, $false, $false, $false, $true, 1, $false, "97-26=71", 2) | Out-Null
$d.Content.Find.Execute("79+4=83", $true, $false, $false, $false, $false, $true, 1, $false, "61-6=55", 2) | Out-Null
$d.Content.Find.Execute("76+2=78", $true, $false, $false, $false, $false, $true, 1, $false, "21+26=47", 2) | Out-Null
$d.Content.Find.Execute("52-38=14", $true, $false, $false, $false, $false, $true, 1, $false, "87-7=80", 2) | Out-Null
$d.Content.Find.Execute("38-8=30", $true, $false, $false, $false, $false, $true, 1, $false, "86-61=25", 2) | Out-Null
$d.Content.Find.Execute("18+45=63", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=67", 2) | Out-Null
$d.Content.Find.Execute("5+30=35", $true, $false, $false, $false, $false, $true, 1, $false, "48+44=92", 2) | Out-Null
$d.Content.Find.Execute("91-82=9", $true, $false, $false, $false, $false, $true, 1, $false, "24+37=61", 2) | Out-Null
$d.Content.Find.Execute("5+35=40", $true, $false, $false, $false, $false, $true, 1, $false, "18+65=83", 2) | Out-Null
$d.Content.Find.Execute("87-83=4", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=76", 2) | Out-Null
$d.Content.Find.Execute("48+17=65", $true, $false, $false, $false, $false, $true, 1, $false, "27+27=54", 2) | Out-Null
$d.Content.Find.Execute("20+37=57", $true, $false, $false, $false, $false, $true, 1, $false, "69+28=97", 2) | Out-Null
$d.Content.Find.Execute("19+33=52", $true, $false, $false, $false, $false, $true, 1, $false, "10+83=93", 2) | Out-Null
$d.Content.Find.Execute("92+5=97", $true, $false, $false, $false, $false, $true, 1, $false, "79-4=75", 2) | Out-Null
$d.Content.Find.Execute("68-28=40", $true, $false, $false, $false, $false, $true, 1, $false, "39+10=49", 2) | Out-Null
$d.Content.Find.Execute("33+16=49", $true, $false, $false, $false, $false, $true, 1, $false, "91-7=84", 2) | Out-Null
$d.Content.Find.Execute("48-7=41", $true, $false, $false, $false, $false, $true, 1, $false, "58-42=16", 2) | Out-Null
$d.Content.Find.Execute("58-40=18", $true, $false, $false, $false, $false, $true, 1, $false, "78+0=78", 2) | Out-Null
$d.Content.Find.Execute("4+41=45", $true, $false, $false, $false, $false, $true, 1, $false, "1+85=86", 2) | Out-Null
$d.Content.Find.Execute("70-12=58", $true, $false, $false, $false, $false, $true, 1, $false, "41+2=43", 2) | Out-Null
$d.Content.Find.Execute("72-32=40", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=89", 2) | Out-Null
$d.Content.Find.Execute("26+24=50", $true, $false, $false, $false, $false, $true, 1, $false, "40-26=14", 2) | Out-Null
$d.Content.Find.Execute("64-30=34", $true, $false, $false, $false, $false, $true, 1, $false, "94-83=11", 2) | Out-Null
$d.Content.Find.Execute("95-82=13", $true, $false, $false, $false, $false, $true, 1, $false, "97-59=38", 2) | Out-Null
$d.Content.Find.Execute("47+24=71", $true, $false, $false, $false, $false, $true, 1, $false, "73-73=0", 2) | Out-Null
$d.Content.Find.Execute("42-36=6", $true, $false, $false, $false, $false, $true, 1, $false, "65-12=53", 2) | Out-Null
$d.Content.Find.Execute("18+11=29", $true, $false, $false, $false, $false, $true, 1, $false, "49+21=70", 2) | Out-Null
$d.Content.Find.Execute("74-52=22", $true, $false, $false, $false, $false, $true, 1, $false, "68+11=79", 2) | Out-Null
$d.Content.Find.Execute("69+29=98", $true, $false, $false, $false, $false, $true, 1, $false, "36-4=32", 2) | Out-Null
$d.Content.Find.Execute("36+35=71", $true, $false, $false, $false, $false, $true, 1, $false, "14+66=80", 2) | Out-Null
$d.Content.Find.Execute("30+13=43", $true, $false, $false, $false, $false, $true, 1, $false, "69+18=87", 2) | Out-Null
$d.Content.Find.Execute("50+7=57", $true, $false, $false, $false, $false, $true, 1, $false, "89-24=65", 2) | Out-Null
$d.Content.Find.Execute("86-81=5", $true, $false, $false, $false, $false, $true, 1, $false, "81-40=41", 2) | Out-Null
$d.Content.Find.Execute("79-22=57", $true, $false, $false, $false, $false, $true, 1, $false, "99-42=57", 2) | Out-Null
$d.Content.Find.Execute("54+32=86", $true, $false, $false, $false, $false, $true, 1, $false, "77-59=18", 2) | Out-Null
$d.Content.Find.Execute("1+49=50", $true, $false, $false, $false, $false, $true, 1, $false, "21+73=94", 2) | Out-Null
$d.Content.Find.Execute("43+51=94", $true, $false, $false, $false, $false, $true, 1, $false, "21+72=93", 2) | Out-Null
$d.Content.Find.Execute("76-24=52", $true, $false, $false, $false, $false, $true, 1, $false, "7+33=40", 2) | Out-Null
$d.Content.Find.Execute("2+78=80", $true, $false, $false, $false, $false, $true, 1, $false, "47+4=51", 2) | Out-Null
$d.Content.Find.Execute("12+9=21", $true, $false, $false, $false, $false, $true, 1, $false, "92-27=65", 2) | Out-Null
$d.Content.Find.Execute("3+47=50", $true, $false, $false, $false, $false, $true, 1, $false, "6+24=30", 2) | Out-Null
$d.Content.Find.Execute("35+62=97", $true, $false, $false, $false, $false, $true, 1, $false, "44-11=33", 2) | Out-Null
$d.Content.Find.Execute("49-13=36", $true, $false, $false, $false, $false, $true, 1, $false, "28-5=23", 2) | Out-Null
$d.Content.Find.Execute("24-11=13", $true, $false, $false, $false, $false, $true, 1, $false, "93-76=17", 2) | Out-Null
$d.Content.Find.Execute("91-66=25", $true, $false, $false, $false, $false, $true, 1, $false, "89-5=84", 2) | Out-Null
$d.Content.Find.Execute("7+44=51", $true, $false, $false, $false, $false, $true, 1, $false, "65-39=26", 2) | Out-Null
$d.Content.Find.Execute("14+46=60", $true, $false, $false, $false, $false, $true, 1, $false, "43-15=28", 2) | Out-Null
